$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.494.26'
$ws.Range("E2").Value = '  +7.73%  '

$ws.Range("D3").Value = '3.508.99'
$ws.Range("E3").Value = '  +9.50%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '189.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +13.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '554.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.04%  '

$ws.Range("D7").Value = '3.501.52'
$ws.Range("E7").Value = '  +9.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.612'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.81%  '

$ws.Range("E9").Value = '  -0.21%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.639'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.98'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.10%  '

$ws.Range("E12").Value = '  +17.03%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000276'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.51'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.76%  '

$ws.Range("D15").Value = '4.063.30'
$ws.Range("E15").Value = '  +8.75%  '

$ws.Range("D16").Value = '3.502.47'
$ws.Range("E16").Value = '  +8.96%  '

$ws.Range("D17").Value = '68.320.31'
$ws.Range("E17").Value = '  +9.09%  '

$ws.Range("E18").Value = '  +7.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.96%  '

$ws.Range("E21").Value = '  +9.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '408.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +13.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.63%  '

$ws.Range("E24").Value = '  +11.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.63'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.39%  '

$ws.Range("E26").Value = '  +11.45%  '

$ws.Range("E27").Value = '  +3.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.91'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '688.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.43%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.84'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.15%  '

$ws.Range("E34").Value = '  +7.73%  '

$ws.Range("E35").Value = '  +9.57%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '60.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '39.18'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.13%  '

$ws.Range("D38").Value = '0.0₃0829'
$ws.Range("E38").Value = '  +25.74%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.402'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.81%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +25.96%  '

$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.133'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.59%  '

$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.78'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +17.94%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +18.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("D46").Value = '3.056.46'
$ws.Range("E46").Value = '  +7.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0422'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.48%  '

$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +24.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.50%  '

$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.25'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.99%  '

$ws.Range("E51").Value = '  +7.78%  '
